$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:50:53"
$ws1.Range("A3").Value = "Total filas: 64"
$ws1.Range("A6:E1000").Clear()

$text1 = @"
03:45:24|03:46|14_ABASTO|1|LP1912
03:45:24|04:01|81_EL PELIGRO|16|LP1912
04:38:41|04:45|215A_EL PATO|7|LP1912
03:45:24|04:46|215A_EL PATO|61|LP1912
04:38:41|04:53|11_ETCHEVERRY|15|LP1912
05:05:17|05:16|17_ROMERO|11|LP1912
04:38:41|05:22|23_HERNANDEZ|44|LP1912
05:05:17|05:23|23_HERNANDEZ|18|LP1912
04:38:41|05:34|215B_EL PATO|56|LP1912
05:05:17|05:35|215B_EL PATO|30|LP1912
03:45:24|05:36|14_ABASTO|111|LP1912
05:05:17|05:46|15_ABASTO|41|LP1912
05:05:17|05:54|10_OLMOS|49|LP1912
05:59:00|06:00|81_EL PELIGRO|1|LP1912
05:05:17|06:04|16_SANTA ANA|59|LP1912
05:05:17|06:11|215A_EL PATO|66|LP1912
05:59:00|06:12|215A_EL PATO|13|LP1912
05:59:00|06:14|225_HARAS DEL SUR|15|LP1912
05:59:00|06:21|26_HERNANDEZ|22|LP1912
05:59:00|06:27|23_HERNANDEZ|28|LP1912
04:38:41|06:29|86_EST CHICA-ESC AGRARIA|111|LP1912
05:59:00|06:30|86_EST CHICA-ESC AGRARIA|31|LP1912
05:05:17|06:31|16_SANTA ANA|86|LP1912
05:59:00|06:32|16_SANTA ANA|33|LP1912
05:59:00|06:44|225_C ROCA-H SUR|45|LP1912
05:05:17|06:46|215C_EL PATO|101|LP1912
05:59:00|06:47|215C_EL PATO|48|LP1912
06:50:53|06:50|215C_EL PATO|0|LP1912
06:50:53|06:59|14_ABASTO|9|LP1912
05:59:00|07:00|14_ABASTO|61|LP1912
06:50:53|07:01|16_SANTA ANA|11|LP1912
06:50:53|07:04|23_HERNANDEZ|14|LP1912
05:59:00|07:05|23_HERNANDEZ|66|LP1912
06:50:53|07:05|15_ABASTO|15|LP1912
06:50:53|07:06|225_GOMEZ|16|LP1912
05:59:00|07:07|225_GOMEZ|68|LP1912
06:50:53|07:11|215A_EL PATO|21|LP1912
05:59:00|07:12|215A_EL PATO|73|LP1912
06:50:53|07:13|16_SANTA ANA|23|LP1912
06:50:53|07:15|11_ETCHEVERRY|25|LP1912
05:59:00|07:16|11_ETCHEVERRY|77|LP1912
06:50:53|07:20|26_HERNANDEZ|30|LP1912
05:59:00|07:21|26_HERNANDEZ|82|LP1912
06:50:53|07:22|10_OLMOS|32|LP1912
05:59:00|07:23|10_OLMOS|84|LP1912
06:50:53|07:27|10_OLMOS|37|LP1912
06:50:53|07:31|11_ETCHEVERRY|41|LP1912
06:50:53|07:32|84_COLONIA URQUIZA-ESC 49|42|LP1912
05:59:00|07:32|16_SANTA ANA|93|LP1912
05:59:00|07:32|11_ETCHEVERRY|93|LP1912
06:50:53|07:36|27_EL RETIRO|46|LP1912
05:59:00|07:37|27_EL RETIRO|98|LP1912
05:59:00|07:39|10_OLMOS|100|LP1912
06:50:53|07:47|14_ABASTO|57|LP1912
05:59:00|07:48|14_ABASTO|109|LP1912
06:50:53|07:51|215D_EL PATO|61|LP1912
05:59:00|07:52|215D_EL PATO|113|LP1912
06:50:53|08:03|23_HERNANDEZ|73|LP1912
06:50:53|08:12|15_ABASTO|82|LP1912
06:50:53|08:20|26_HERNANDEZ|90|LP1912
06:50:53|08:22|16_P MOR-SANTA ANA|92|LP1912
06:50:53|08:23|215B_EL PATO|93|LP1912
06:50:53|08:27|84_COLONIA URQUIZA-ESC 49|97|LP1912
06:50:53|08:42|81_EL PELIGRO|112|LP1912
"@
$lines1 = $text1 -split "`n"
$r = 6
foreach ($l in $lines1) {
  $l = $l.Trim()
  if ($l.Length -eq 0) { continue }
  $parts = $l -split '\|'
  $ws1.Cells.Item($r,1).Value = $parts[0]
  $ws1.Cells.Item($r,2).Value = $parts[1]
  $ws1.Cells.Item($r,3).Value = $parts[2]
  $ws1.Cells.Item($r,4).Value = [int]$parts[3]
  $ws1.Cells.Item($r,5).Value = $parts[4]
  $r = $r + 1
}

# ---------- Sheet 2: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:50:53"
$ws2.Range("A3").Value = "Total filas: 14"
$ws2.Range("A6:E1000").Clear()

$text2 = @"
04:38:41|04:45|215A_EL PATO|7|LP1912
03:45:24|04:46|215A_EL PATO|61|LP1912
04:38:41|05:34|215B_EL PATO|56|LP1912
05:05:17|05:35|215B_EL PATO|30|LP1912
05:05:17|06:11|215A_EL PATO|66|LP1912
05:59:00|06:12|215A_EL PATO|13|LP1912
05:05:17|06:46|215C_EL PATO|101|LP1912
05:59:00|06:47|215C_EL PATO|48|LP1912
06:50:53|06:50|215C_EL PATO|0|LP1912
06:50:53|07:11|215A_EL PATO|21|LP1912
05:59:00|07:12|215A_EL PATO|73|LP1912
06:50:53|07:51|215D_EL PATO|61|LP1912
05:59:00|07:52|215D_EL PATO|113|LP1912
06:50:53|08:23|215B_EL PATO|93|LP1912
"@
$lines2 = $text2 -split "`n"
$r = 6
foreach ($l in $lines2) {
  $l = $l.Trim()
  if ($l.Length -eq 0) { continue }
  $parts = $l -split '\|'
  $ws2.Cells.Item($r,1).Value = $parts[0]
  $ws2.Cells.Item($r,2).Value = $parts[1]
  $ws2.Cells.Item($r,3).Value = $parts[2]
  $ws2.Cells.Item($r,4).Value = [int]$parts[3]
  $ws2.Cells.Item($r,5).Value = $parts[4]
  $r = $r + 1
}

# ---------- Sheet 3: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:50:53"
$ws3.Range("A3").Value = "Total filas: 11"
$ws3.Range("A6:E1000").Clear()

$text3 = @"
04:38:41|05:43|215A_LA PLATA|65|L6173
05:05:17|05:44|215A_LA PLATA|39|L6173
04:38:41|06:08|215A_LA PLATA|90|L6173
05:59:00|06:09|215A_LA PLATA|10|L6173
04:38:41|06:32|215C_LA PLATA|114|L6203
05:59:00|06:33|215C_LA PLATA|34|L6203
06:50:53|06:59|215B_LP-P MOR-1 Y 57|9|L6173
05:59:00|07:00|215B_LP-P MOR-1 Y 57|61|L6173
06:50:53|07:35|215A_LA PLATA|45|L6173
06:50:53|08:08|215C_LA PLATA|78|L6203
06:50:53|08:35|215A_LA PLATA|105|L6173
"@
$lines3 = $text3 -split "`n"
$r = 6
foreach ($l in $lines3) {
  $l = $l.Trim()
  if ($l.Length -eq 0) { continue }
  $parts = $l -split '\|'
  $ws3.Cells.Item($r,1).Value = $parts[0]
  $ws3.Cells.Item($r,2).Value = $parts[1]
  $ws3.Cells.Item($r,3).Value = $parts[2]
  $ws3.Cells.Item($r,4).Value = [int]$parts[3]
  $ws3.Cells.Item($r,5).Value = $parts[4]
  $r = $r + 1
}

Write-Output "Done"
